$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update end time of the last working entry (row 42) from 22:00 to 22:30.
$ws.Range("E42").Value = 0.9375

# Update the active selection to match the saved view (E43).
$ws.Range("E43").Select()
